# Rename the "gaia_id" header (column A) to "source_id"
# (the workbook now also accommodates Bailer-Jones distance IDs, not just Gaia IDs)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "source_id"

# Reset the selection / active cell back to A1 so the sheet does not keep
# the stray "D17" selection that was stored in the file before editing.
$ws.Range("A1").Select()
